# Add season-record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (A1) onto the new headers
# so they share the same style as the rest of row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2 through 48 get the season record values
$wins = 84
$losses = 78
$ties = 0

for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($r, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($r, 32).Value = $ties    # column AF = 32
}
